# FLRenamingConfigFile.xlsx - Update ConfigFile (To and CC list) for testing in Orchestrator.
#
# Summary of the data change:
#  - "Config" sheet: the RecipientTo ("To List") and RecipientCC ("CC List") values
#    are both replaced with a single address: lester.rollan@lexisnexisrisk.com
#    (rows shrink back to default height since the new text is short).
#  - "ConfigOptions" sheet (3 repeated config blocks): the RecipientCC ("CC List")
#    value gains two more recipients (paul.fabro, judy.cotaoco) - the RecipientTo
#    ("To List") value is unchanged. The CC row grows taller to fit the longer text.
#  - The view/selection on "ConfigOptions" is moved down toward the bottom block.

$wb = $excel.ActiveWorkbook

$newSingleAddress = "lester.rollan@lexisnexisrisk.com"
$newCcList = "lester.rollan@lexisnexisrisk.com; dindee.galindo@lexisnexisrisk.com; jesriel.tolentino@lexisnexisrisk.com; jhoanna.talle@lexisnexisrisk.com; paul.fabro@lexisnexisrisk.com; judy.cotaoco@lexisnexisrisk.com"

# ---------------------------------------------------------------------------
# 1) "Config" sheet - RecipientTo / RecipientCC rows (10 and 11)
# ---------------------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("Config")

$wsConfig.Range("B10").Value = $newSingleAddress
$wsConfig.Range("B11").Value = $newSingleAddress

# Rows no longer need the taller wrapped height - let them size back down.
$wsConfig.Rows.Item(10).AutoFit()
$wsConfig.Rows.Item(11).AutoFit()

# ---------------------------------------------------------------------------
# 2) "ConfigOptions" sheet - three repeated config blocks
#    (rows 11/12, 25/26, 39/40) - only the CC List value + its row height change.
# ---------------------------------------------------------------------------
$wsOptions = $wb.Worksheets.Item("ConfigOptions")

$wsOptions.Range("B12").Value = $newCcList
$wsOptions.Rows.Item(12).RowHeight = 45

$wsOptions.Range("B26").Value = $newCcList
$wsOptions.Rows.Item(26).RowHeight = 45

$wsOptions.Range("B40").Value = $newCcList
$wsOptions.Rows.Item(40).RowHeight = 45

# ---------------------------------------------------------------------------
# 3) Update the view/selection on "ConfigOptions" to the bottom block, then
#    restore "Config" as the active sheet/tab (matches original workbook state).
# ---------------------------------------------------------------------------
$wsOptions.Activate()
$wsOptions.Range("B35:C38").Select()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1

$wsConfig.Activate()
$wsConfig.Range("A1").Select()
